# user story ## 20 Completed ## 1 hour, an email is now able to be sent with
# a link for a user to reset their password.
#
# The rubric row for that story ("Forgot password" page that emails reset
# links, with exp date -> B23) had no "Points" score recorded yet (C23 was
# blank). Now that the story is complete, record the 10 points earned.
# The Points-possible (D23) already reads 10, so this maxes out that row.
# All of the roll-up totals (E18, E25, E27, F28) are formulas and
# recalculate automatically from this single input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the 10 points scored for the "Forgot password" reset-link story.
$ws.Range("C23").Value = 10

# Reflect where the grader's cursor/viewport ended up after making the edit.
$ws.Range("C24").Select()
